$d = $word.ActiveDocument

# Locate the paragraph containing "Report for Coffee Plant information".
# It is immediately followed by two empty paragraphs and then the
# "Report sa Sources : " paragraph (which currently carries the
# hidden _GoBack bookmark at its end).
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Report for Coffee Plant information*") {
        $anchorIndex = $i
        break
    }
}

$emptyIndex1 = $anchorIndex + 1
$emptyIndex2 = $anchorIndex + 2

# Collapse the two empty paragraphs into one by deleting the second
# empty paragraph's range (this removes it along with its paragraph mark).
$d.Paragraphs.Item($emptyIndex2).Range.Delete()

# Move the _GoBack bookmark off of the "Report sa Sources : " paragraph
# and onto the now-single remaining empty paragraph above it.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$remainingEmpty = $d.Paragraphs.Item($emptyIndex1)
$d.Bookmarks.Add("_GoBack", $remainingEmpty.Range)
